$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of key names (case-only corrections)
$ws.Range("A2").Value = "mdaTextHomePage"
$ws.Range("A8").Value = "pageTitleNewTab"
$ws.Range("A4").Value = "mdaTitle"

# Update the active selection to B9
$ws.Range("B9").Select()
